$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header in H1, matching the style of the other header cells (B1:G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding value in H2
$ws.Range("H2").Value = 1
